$wb = $excel.ActiveWorkbook
Write-Output $wb.Worksheets.Count
